$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PSFBeadsKeyMeasurements")

# Insert a new column at C, shifting considered_valid_count etc. to the right
$ws.Columns.Item(3).Insert()

# Set the header for the newly inserted column
$ws.Cells.Item(1, 3).Value = "total_bead_count"
